# Daily attendance processing - reorder the "Recorded By" (column G) list
# for rows whose value contains the literal token "System" (but is not an
# admin@admin.com entry) by reversing the comma-separated author list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if ($value -notmatch ",") { continue }

    $parts = @($value -split ",\s*" | ForEach-Object { $_.Trim() })

    if (($parts -contains "System") -and ($parts -notcontains "admin@admin.com")) {
        $reversed = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $cell.Value2 = [string]::Join(", ", $reversed)
    }
}
